$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: add P1=14, Q1=15 (continuing the 0..13 sequence), same style as existing header row (s=1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Rows 2-25: swap values in columns I, K, M, O and add new columns P, Q (=2)
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new
    $ws.Cells.Item($r, 17).Value = 2   # Q: new
}
